$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C8) from 2023-10-09 (45208) to 2023-10-13 (45212)
$ws.Range("C2:C8").Value = 45212
